# Auto-generated edit script: bump ObjTables schema metadata (date/version)
# and rename the Relation/Relationship "From"/"To" fields to "FromObject"/"ToObject".
$wb = $excel.ActiveWorkbook

# --- 1. Update ObjTables header metadata (date + objTablesVersion) on every sheet ---
$ws___Compartment = $wb.Worksheets.Item("!!Compartment")
$ws___Compartment.Unprotect()
$ws___Compartment.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.9' date='2020-04-26 21:08:05'"
$ws___Compartment.Range("A2").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compartment' name='Compartment' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Compartment.Protect()

$ws___Compound = $wb.Worksheets.Item("!!Compound")
$ws___Compound.Unprotect()
$ws___Compound.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compound' name='Compound' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Compound.Protect()

$ws___Definition = $wb.Worksheets.Item("!!Definition")
$ws___Definition.Unprotect()
$ws___Definition.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Definition' name='Definition' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Definition.Protect()

$ws___Enzyme = $wb.Worksheets.Item("!!Enzyme")
$ws___Enzyme.Unprotect()
$ws___Enzyme.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Enzyme' name='Enzyme' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Enzyme.Protect()

$ws___FbcObjective = $wb.Worksheets.Item("!!FbcObjective")
$ws___FbcObjective.Unprotect()
$ws___FbcObjective.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='FbcObjective' name='FbcObjective' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___FbcObjective.Protect()

$ws___Gene = $wb.Worksheets.Item("!!Gene")
$ws___Gene.Unprotect()
$ws___Gene.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Gene' name='Gene' date='2020-04-26 21:08:05' objTablesVersion='0.0.9' document='lac_Operon'"
$ws___Gene.Protect()

$ws___Layout = $wb.Worksheets.Item("!!Layout")
$ws___Layout.Unprotect()
$ws___Layout.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Layout' name='Layout' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Layout.Protect()

$ws___Measurement = $wb.Worksheets.Item("!!Measurement")
$ws___Measurement.Unprotect()
$ws___Measurement.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Measurement' name='Measurement' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Measurement.Protect()

$ws___PbConfig = $wb.Worksheets.Item("!!PbConfig")
$ws___PbConfig.Unprotect()
$ws___PbConfig.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='PbConfig' name='PbConfig' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___PbConfig.Protect()

$ws___Position = $wb.Worksheets.Item("!!Position")
$ws___Position.Unprotect()
$ws___Position.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Position' name='Position' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Position.Protect()

$ws___Protein = $wb.Worksheets.Item("!!Protein")
$ws___Protein.Unprotect()
$ws___Protein.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Protein' name='Protein' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Protein.Protect()

$ws___Quantity = $wb.Worksheets.Item("!!Quantity")
$ws___Quantity.Unprotect()
$ws___Quantity.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Quantity' name='Quantity' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Quantity.Protect()

$ws___QuantityInfo = $wb.Worksheets.Item("!!QuantityInfo")
$ws___QuantityInfo.Unprotect()
$ws___QuantityInfo.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityInfo' name='QuantityInfo' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___QuantityInfo.Protect()

$ws___QuantityMatrix = $wb.Worksheets.Item("!!QuantityMatrix")
$ws___QuantityMatrix.Unprotect()
$ws___QuantityMatrix.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityMatrix' name='QuantityMatrix' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___QuantityMatrix.Protect()

$ws___Reaction = $wb.Worksheets.Item("!!Reaction")
$ws___Reaction.Unprotect()
$ws___Reaction.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Reaction' name='Reaction' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Reaction.Protect()

$ws___ReactionStoichiometry = $wb.Worksheets.Item("!!ReactionStoichiometry")
$ws___ReactionStoichiometry.Unprotect()
$ws___ReactionStoichiometry.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___ReactionStoichiometry.Protect()

$ws___Regulator = $wb.Worksheets.Item("!!Regulator")
$ws___Regulator.Unprotect()
$ws___Regulator.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Regulator' name='Regulator' date='2020-04-26 21:08:05' objTablesVersion='0.0.9' document='lac_Operon'"
$ws___Regulator.Protect()

$ws___Relation = $wb.Worksheets.Item("!!Relation")
$ws___Relation.Unprotect()
$ws___Relation.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relation' name='Relation' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Relation.Protect()

$ws___Relationship = $wb.Worksheets.Item("!!Relationship")
$ws___Relationship.Unprotect()
$ws___Relationship.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relationship' name='Relationship' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___Relationship.Protect()

$ws___SparseMatrix = $wb.Worksheets.Item("!!SparseMatrix")
$ws___SparseMatrix.Unprotect()
$ws___SparseMatrix.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrix' name='SparseMatrix' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___SparseMatrix.Protect()

$ws___SparseMatrixColumn = $wb.Worksheets.Item("!!SparseMatrixColumn")
$ws___SparseMatrixColumn.Unprotect()
$ws___SparseMatrixColumn.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___SparseMatrixColumn.Protect()

$ws___SparseMatrixOrdered = $wb.Worksheets.Item("!!SparseMatrixOrdered")
$ws___SparseMatrixOrdered.Unprotect()
$ws___SparseMatrixOrdered.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___SparseMatrixOrdered.Protect()

$ws___SparseMatrixRow = $wb.Worksheets.Item("!!SparseMatrixRow")
$ws___SparseMatrixRow.Unprotect()
$ws___SparseMatrixRow.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixRow' name='SparseMatrixRow' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___SparseMatrixRow.Protect()

$ws___StoichiometricMatrix = $wb.Worksheets.Item("!!StoichiometricMatrix")
$ws___StoichiometricMatrix.Unprotect()
$ws___StoichiometricMatrix.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___StoichiometricMatrix.Protect()

$ws___rxnconContingencyList = $wb.Worksheets.Item("!!rxnconContingencyList")
$ws___rxnconContingencyList.Unprotect()
$ws___rxnconContingencyList.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconContingencyList' name='rxnconContingencyList' date='2020-04-26 21:08:05' objTablesVersion='0.0.9'"
$ws___rxnconContingencyList.Protect()

$ws___rxnconReactionList = $wb.Worksheets.Item("!!rxnconReactionList")
$ws___rxnconReactionList.Unprotect()
$ws___rxnconReactionList.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconReactionList' name='rxnconReactionList' date='2020-04-26 21:08:06' objTablesVersion='0.0.9'"
$ws___rxnconReactionList.Protect()

# --- 2. Rename "From"/"To" columns to "FromObject"/"ToObject" on !!Relation and !!Relationship ---
$wsRelation = $wb.Worksheets.Item("!!Relation")
$wsRelation.Unprotect()
$wsRelation.Range("G2").Value = "!FromObject"
$wsRelation.Range("H2").Value = "!ToObject"

$valG = $wsRelation.Range("G2").Validation
$valG.ErrorTitle = "FromObject"
$valG.InputTitle = "FromObject"

$valH = $wsRelation.Range("H2").Validation
$valH.ErrorTitle = "ToObject"
$valH.InputTitle = "ToObject"
$wsRelation.Protect()

$wsRelationship = $wb.Worksheets.Item("!!Relationship")
$wsRelationship.Unprotect()
$wsRelationship.Range("B2").Value = "!FromObject"
$wsRelationship.Range("C2").Value = "!ToObject"

$valB = $wsRelationship.Range("B2").Validation
$valB.ErrorTitle = "FromObject"
$valB.InputTitle = "FromObject"

$valC = $wsRelationship.Range("C2").Validation
$valC.ErrorTitle = "ToObject"
$valC.InputTitle = "ToObject"
$wsRelationship.Protect()
